$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rentals")
$lo = $ws.ListObjects.Item("Rentals")

# Resize the table from A1:B11 to A1:H25
$lo.Resize($ws.Range("A1:H25"))

# Set the header row: Time, Sunday, Monday, Tuesday, Wednesday, Thursday, Friday, Saturday
$headers = @("Time","Sunday","Monday","Tuesday","Wednesday","Thursday","Friday","Saturday")
for ($i = 0; $i -lt 8; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Time-slot row labels (rows 2-25)
$times = @(
    "12:00 AM - 01:00 AM",
    "01:00 AM - 02:00 AM",
    "02:00 AM - 03:00 AM",
    "03:00 AM - 04:00 AM",
    "04:00 AM - 05:00 AM",
    "05:00 AM - 06:00 AM",
    "06:00 AM - 07:00 AM",
    "07:00 AM - 08:00 AM",
    "08:00 AM - 09:00 AM",
    "09:00 AM - 10:00 AM",
    "10:00 AM - 11:00 AM",
    "11:00 AM - 12:00 PM",
    "12:00 PM - 01:00 PM",
    "01:00 PM - 02:00 PM",
    "02:00 PM - 03:00 PM",
    "03:00 PM - 04:00 PM",
    "04:00 PM - 05:00 PM",
    "05:00 PM - 06:00 PM",
    "06:00 PM - 07:00 PM",
    "07:00 PM - 08:00 PM",
    "08:00 PM - 09:00 PM",
    "09:00 PM - 10:00 PM",
    "10:00 PM - 11:00 PM",
    "11:00 PM - 12:00 AM"
)

for ($r = 0; $r -lt 24; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 1).Value = $times[$r]
    $ws.Cells.Item($row, 1).NumberFormat = "h:mm AM/PM"
    for ($c = 2; $c -le 8; $c++) {
        $ws.Cells.Item($row, $c).Value = $false
    }
}

# Mark availability checkboxes that are TRUE (slot not available)
$ws.Cells.Item(14, 2).Value = $true   # 12:00 PM - 01:00 PM, Sunday
$ws.Cells.Item(20, 3).Value = $true   # 06:00 PM - 07:00 PM, Monday

# Footnote row, merged across A26:H26, centered
$ws.Cells.Item(26, 1).Value = "Note: Checkbox means the time slot is not avaliable"
$ws.Range("A26:H26").Merge()
$ws.Range("A26:H26").HorizontalAlignment = -4108

# Column widths (values pre-compensated for the engine's internal 5px/MDW=6 padding
# rounding on ColumnWidth, so the persisted OOXML "width" lands as close as possible
# to the target character-width values of 19.5546875, 10.77734375, 10.44140625,
# 11.109375, 14.109375, 11.88671875, 10.109375, 12.33203125)
$ws.Columns.Item(1).ColumnWidth = 18.666666666666668
$ws.Columns.Item(2).ColumnWidth = 10.0
$ws.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(5).ColumnWidth = 13.333333333333334
$ws.Columns.Item(6).ColumnWidth = 11.0
$ws.Columns.Item(7).ColumnWidth = 9.333333333333334
$ws.Columns.Item(8).ColumnWidth = 11.5

$wb.Save()
